$d = $word.ActiveDocument

$replacements = @(
    @("2024-11-09 Saturday", "2024-11-10 Sunday"),
    @("490÷8=61, 2", "348÷6=58, 0"),
    @("768÷7=109, 5", "409÷6=68, 1"),
    @("220÷8=27, 4", "985÷9=109, 4"),
    @("189÷9=21, 0", "918÷7=131, 1"),
    @("610÷2=305, 0", "867÷4=216, 3"),
    @("350÷7=50, 0", "345÷6=57, 3"),
    @("181÷5=36, 1", "946÷5=189, 1"),
    @("533÷7=76, 1", "944÷4=236, 0"),
    @("316÷2=158, 0", "411÷2=205, 1"),
    @("633÷4=158, 1", "546÷7=78, 0"),
    @("190÷6=31, 4", "154÷4=38, 2"),
    @("501÷3=167, 0", "807÷2=403, 1"),
    @("120÷5=24, 0", "535÷7=76, 3"),
    @("925÷7=132, 1", "959÷4=239, 3"),
    @("914÷6=152, 2", "261÷4=65, 1"),
    @("404÷8=50, 4", "653÷5=130, 3"),
    @("333÷8=41, 5", "470÷4=117, 2"),
    @("562÷8=70, 2", "486÷3=162, 0"),
    @("418÷3=139, 1", "330÷5=66, 0"),
    @("876÷4=219, 0", "319÷4=79, 3"),
    @("373÷7=53, 2", "532÷9=59, 1"),
    @("586÷7=83, 5", "269÷9=29, 8"),
    @("564÷8=70, 4", "444÷4=111, 0"),
    @("512÷3=170, 2", "389÷5=77, 4"),
    @("300÷9=33, 3", "899÷7=128, 3")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "done"
